$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29-137 down to 30-138.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new data record.
$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(29, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(29, 4).Value = 45145
$ws.Cells.Item(29, 5).Value = 15
$ws.Cells.Item(29, 6).Value = 100112038
$ws.Cells.Item(29, 7).Value = "Cebollín baby"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 250
$ws.Cells.Item(29, 11).Value = 2800
$ws.Cells.Item(29, 12).Value = 3000
$ws.Cells.Item(29, 13).Value = 2920
$ws.Cells.Item(29, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(29, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(29, 16).Value = 1460
$ws.Cells.Item(29, 17).Value = 2
$ws.Cells.Item(29, 18).Value = "Hortaliza"
